$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.525.30'
$ws.Range("E2").Value = '  +3.90%  '
$ws.Range("D3").Value = '''1.738.60'
$ws.Range("D4").Value = '''0.9997'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''245.69'
$ws.Range("E5").Value = '  +4.83%  '
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = '''0.4804'
$ws.Range("E7").Value = '  +3.36%  '
$ws.Range("E8").Value = '  +4.06%  '
$ws.Range("D9").Value = '''0.06243'
$ws.Range("E9").Value = '  +1.81%  '
$ws.Range("D10").Value = '''1.738.26'
$ws.Range("E10").Value = '  +4.38%  '
$ws.Range("D11").Value = '''0.07128'
$ws.Range("E11").Value = '  +2.65%  '
$ws.Range("E12").Value = '  +7.95%  '
$ws.Range("D13").Value = '''0.6208'
$ws.Range("E13").Value = '  +8.54%  '
$ws.Range("D14").Value = '''4.545'
$ws.Range("E14").Value = '  +4.30%  '
$ws.Range("D15").Value = '''77.17'
$ws.Range("E15").Value = '  +3.00%  '
$ws.Range("D16").Value = '''1.000'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = '''26.533.48'
$ws.Range("E17").Value = '  +3.93%  '
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '''0.000006893'
$ws.Range("E19").Value = '  +2.71%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '''11.77'
$ws.Range("E20").Value = '  +3.64%  '
$ws.Range("D21").Value = '''1.961.52'
$ws.Range("E21").Value = '  +4.43%  '
$ws.Range("D22").Value = '''4.601'
$ws.Range("E22").Value = '  +4.11%  '
$ws.Range("D23").Value = '''8.901'
$ws.Range("E23").Value = '  +2.39%  '
$ws.Range("E24").Value = '  +2.69%  '
$ws.Range("D25").Value = '''136.22'
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("D26").Value = '''15.39'
$ws.Range("E26").Value = '  +3.46%  '
$ws.Range("D27").Value = '''1.809'
$ws.Range("E27").Value = '  +6.16%  '
$ws.Range("D28").Value = '''1.427'
$ws.Range("E28").Value = '  +4.64%  '
$ws.Range("D29").Value = '''106.93'
$ws.Range("E29").Value = '  +2.94%  '
$ws.Range("D30").Value = '''3.994'
$ws.Range("E30").Value = '  +1.20%  '
$ws.Range("D31").Value = '''3.740'
$ws.Range("E31").Value = '  +3.88%  '
$ws.Range("D32").Value = '''0.07874'
$ws.Range("E32").Value = '  +2.22%  '
$ws.Range("D33").Value = '''0.04588'
$ws.Range("E33").Value = '  +6.63%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '''2.613'
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '''1.000'
$ws.Range("E35").Value = '  +6.06%  '
$ws.Range("D36").Value = '''0.6372'
$ws.Range("E36").Value = '  +6.54%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '''0.9292'
$ws.Range("E37").Value = '  +0.91%  '
$ws.Range("B38").Value = 'Quant'
$ws.Range("C38").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D38").Value = '''112.89'
$ws.Range("E38").Value = '  +9.42%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '''2.435'
$ws.Range("E39").Value = '  -1.84%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '''1.982'
$ws.Range("E40").Value = '  +8.71%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '''1.005'
$ws.Range("E41").Value = '  +0.59%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.01517'
$ws.Range("E42").Value = '  +3.81%  '
$ws.Range("D43").Value = '''5.736'
$ws.Range("E43").Value = '  +13.01%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '''0.3914'
$ws.Range("E44").Value = '  +5.62%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '''6.931'
$ws.Range("E45").Value = '  +13.41%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '''0.1200'
$ws.Range("E46").Value = '  +8.35%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.05330'
$ws.Range("E47").Value = '  +1.53%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''7.878'
$ws.Range("E48").Value = '  +6.21%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''30.81'
$ws.Range("E49").Value = '  +3.69%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '''1.254'
$ws.Range("E50").Value = '  +5.33%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '''0.3450'
$ws.Range("E51").Value = '  +4.56%  '
